# Apply the commit's changes to the workbook:
#  1. Rename the "data table methods" sheet to "data_table_methods".
#  2. Populate that sheet with the two method rows (title / description /
#     instrumentation) that mirror the existing "methods" sheet, so the
#     template ships with example data for every input sheet.
#  3. Leave the active selection on C3, matching the author's final state.

$wb = $excel.ActiveWorkbook

# 1. Rename the worksheet (sheetId / rId stay the same, only the name changes).
$ws = $wb.Worksheets.Item("data table methods")
$ws.Name = "data_table_methods"

# 2. Add the example rows below the existing header row (title, description, instrumentation).
$ws.Range("A2").Value = "method 1 "
$ws.Range("B2").Value = "this is the first method"
$ws.Range("C2").Value = "ruler"

$ws.Range("A3").Value = "method 2"
$ws.Range("B3").Value = "this is the second method"
$ws.Range("C3").Value = "scale"

# 3. Match the final cell selection recorded in the workbook.
$ws.Range("C3").Select()
